# arrumei o slide 7 do ppt
# The 4th animation effect (fade-in on the big icon picture, shape id 8) was
# starting "On Click" as its own separate click-group; it should instead play
# "With Previous", alongside the other three entrance effects on this slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$seq = $s.TimeLine.MainSequence

$effect = $seq.Item($seq.Count)
$effect.Timing.TriggerType = 2
